$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("M2").Value = 18.79647682669315
$ws.Range("T2").Value = 0.03796096120627238
$ws.Range("V2").Value = 0.01585118999772539
$ws.Range("M3").Value = 18.86992414157695
$ws.Range("T3").Value = 0.03073058057489463
$ws.Range("V3").Value = 0.01295514424609528
$ws.Range("M4").Value = 13.86975728148907
$ws.Range("T4").Value = 0.03690100507967421
$ws.Range("V4").Value = 0.02695082003825436
$ws.Range("M6").Value = 12.21831749797054
$ws.Range("T6").Value = 0.04426490683511781
$ws.Range("V6").Value = 0.02921086162158479
$ws.Range("M10").Value = 13.29779390534793
$ws.Range("T10").Value = 0.03953058041752208
$ws.Range("V10").Value = 0.02903155302235026
$ws.Range("M12").Value = 18.32137548370073
$ws.Range("T12").Value = 0.03685852419268531
$ws.Range("V12").Value = 0.01953613012857874
$ws.Range("M13").Value = 15.98074355023815
$ws.Range("T13").Value = 0.04073603532215879
$ws.Range("V13").Value = 0.02463412059593207
$ws.Range("M14").Value = 19.11910259045317
$ws.Range("T14").Value = 0.02983049032437017
$ws.Range("V14").Value = 0.01696073283616569
$ws.Range("M15").Value = 18.55467266299319
$ws.Range("T15").Value = 0.03023098118544442
$ws.Range("V15").Value = 0.01778316284633774
$ws.Range("M17").Value = 20.63667800776629
$ws.Range("T17").Value = 0.02609499147407053
$ws.Range("V17").Value = 0.01390384055552749
$ws.Range("M18").Value = 12.20975947074733
$ws.Range("M20").Value = 15.53775550141747
$ws.Range("T20").Value = 0.03085124374101027
$ws.Range("V20").Value = 0.02196026262569532
$ws.Range("T21").Value = 0.03863450676709624
$ws.Range("V21").Value = 0.02827116541266036
$ws.Range("M22").Value = 12.39454812261251
$ws.Range("T22").Value = 0.03053860741358308
$ws.Range("V22").Value = 0.02214407354403367
$ws.Range("M24").Value = 11.69632435012107
$ws.Range("T24").Value = 0.03680887107709471
$ws.Range("V24").Value = 0.02935783987608479
$ws.Range("M25").Value = 15.00797764834238
$ws.Range("M26").Value = 20.38148947420772
$ws.Range("M27").Value = 20.33329942796322
$ws.Range("M29").Value = 15.64430077370058
$ws.Range("T29").Value = 0.03454680351546196
$ws.Range("V29").Value = 0.02324392599335197
$ws.Range("M30").Value = 15.63673899668767
$ws.Range("T30").Value = 0.03456384265620741
$ws.Range("V30").Value = 0.02326096513409742
$ws.Range("M31").Value = 13.16126706340616
$ws.Range("T31").Value = 0.04424690196291614
$ws.Range("V31").Value = 0.03112426042970285
$ws.Range("M32").Value = 15.23464679713052
$ws.Range("T32").Value = 0.03649442669593517
$ws.Range("V32").Value = 0.01880167355940041
$ws.Range("M34").Value = 15.51225482349139
$ws.Range("M35").Value = 17.44975671545913
$ws.Range("T35").Value = 0.04551620366825616
$ws.Range("V35").Value = 0.03093550516493822
$ws.Range("M36").Value = 20.29266290898792
$ws.Range("T36").Value = 0.03451904657142538
$ws.Range("V36").Value = 0.01332857535578807
$ws.Range("M37").Value = 20.64144546032226
$ws.Range("T38").Value = 0.02662007470970774
$ws.Range("V38").Value = 0.0156671522231543
$ws.Range("M40").Value = 18.86798492370013
$ws.Range("T40").Value = 0.02939071005607644
$ws.Range("V40").Value = 0.01734715984225846
$ws.Range("M41").Value = 18.75509705334182
$ws.Range("T41").Value = 0.0315959004088202
$ws.Range("V41").Value = 0.01784314606759582
$ws.Range("M42").Value = 12.22462575870084
$ws.Range("T42").Value = 0.03821971878725176
$ws.Range("V42").Value = 0.0237282131256938
$ws.Range("M43").Value = 12.41007311177945
$ws.Range("T43").Value = 0.05131838027354877
$ws.Range("V43").Value = 0.03370308000006633
$ws.Range("M44").Value = 12.31987666854902
$ws.Range("T44").Value = 0.04882314477410078
$ws.Range("V44").Value = 0.02812387631907114
$ws.Range("M45").Value = 16.94624969050644
$ws.Range("T45").Value = 0.03813908253154627
$ws.Range("V45").Value = 0.01684179534228456
$ws.Range("M46").Value = 14.86183663259694
$ws.Range("M47").Value = 15.35194683664891
$ws.Range("M48").Value = 19.67525354469729
$ws.Range("T48").Value = 0.0352049058171427
$ws.Range("V48").Value = 0.01364819352291468
$ws.Range("M67").Value = 12.54450404243209
$ws.Range("T68").Value = 0.04330909613574976
$ws.Range("V68").Value = 0.03369430782764369
$ws.Range("M71").Value = 17.18295849059495
